$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 11, shifting existing rows 11..92 down to 12..93.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with its data.
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(11, 3).Value = "La Araucanía"
$ws.Cells.Item(11, 4).Value = 44881
$ws.Cells.Item(11, 5).Value = 9
$ws.Cells.Item(11, 6).Value = 300000001
$ws.Cells.Item(11, 7).Value = "Rabanito"
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 30
$ws.Cells.Item(11, 11).Value = 8000
$ws.Cells.Item(11, 12).Value = 8000
$ws.Cells.Item(11, 13).Value = 8000
$ws.Cells.Item(11, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(11, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(11, 16).Value = 667
$ws.Cells.Item(11, 17).Value = 12
$ws.Cells.Item(11, 18).Value = "Hortaliza"
